$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target cells whose textual values change (prices / % volumes).
# Force Text format first so values like "1.014" or "0.4821" are kept
# as literal strings instead of being auto-parsed as numbers.
$targetRefs = @(
    'D2', 'E2', 'D3', 'E3', 'D4', 'E4', 'D5', 'E5', 'D6', 'E6',
    'D7', 'E7', 'D8', 'E8', 'D9', 'E9', 'E10', 'D11', 'D12', 'E12',
    'D13', 'E13', 'D14', 'E14', 'D15', 'E15', 'D16', 'E16', 'E17', 'D18',
    'E18', 'E19', 'D20', 'E20', 'D21', 'E21', 'D22', 'E22', 'D23', 'E23',
    'D24', 'E24', 'D25', 'E25', 'D26', 'E26', 'D27', 'E27', 'D28', 'E28',
    'D29', 'E29', 'D30', 'E30', 'D31', 'E31', 'D32', 'E32', 'D33', 'E33',
    'D34', 'E34', 'D35', 'E35', 'D36', 'E36', 'E37', 'D38', 'E38', 'E39',
    'D40', 'E40', 'D41', 'E41', 'D42', 'E42', 'D43', 'E43', 'D44', 'E44',
    'D45', 'E45', 'D46', 'E46', 'D47', 'E47', 'D48', 'E48', 'E49', 'D50',
    'E50', 'D51', 'E51'
)

foreach ($ref in $targetRefs) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range('D2').Value = '29.562.08'
$ws.Range('E2').Value = '  +0.41%  '
$ws.Range('D3').Value = '1.926.86'
$ws.Range('E3').Value = '  +0.95%  '
$ws.Range('D4').Value = '1.014'
$ws.Range('E4').Value = '  +0.78%  '
$ws.Range('D5').Value = '326.65'
$ws.Range('E5').Value = '  +0.40%  '
$ws.Range('D6').Value = '1.013'
$ws.Range('E6').Value = '  +0.69%  '
$ws.Range('D7').Value = '0.4821'
$ws.Range('E7').Value = '  +0.14%  '
$ws.Range('D8').Value = '0.4055'
$ws.Range('E8').Value = '  -0.31%  '
$ws.Range('D9').Value = '0.08203'
$ws.Range('E9').Value = '  +0.53%  '
$ws.Range('E10').Value = '  -0.11%  '
$ws.Range('D11').Value = '23.83'
$ws.Range('D12').Value = '1.947.11'
$ws.Range('E12').Value = '  +3.14%  '
$ws.Range('D13').Value = '6.100'
$ws.Range('E13').Value = '  +1.56%  '
$ws.Range('D14').Value = '7.314'
$ws.Range('E14').Value = '  +2.09%  '
$ws.Range('D15').Value = '91.58'
$ws.Range('E15').Value = '  +1.47%  '
$ws.Range('D16').Value = '0.06896'
$ws.Range('E16').Value = '  +1.89%  '
$ws.Range('E17').Value = '  +0.62%  '
$ws.Range('D18').Value = '0.00001040'
$ws.Range('E18').Value = '  +0.35%  '
$ws.Range('E19').Value = '  -0.03%  '
$ws.Range('D20').Value = '1.012'
$ws.Range('E20').Value = '  +0.66%  '
$ws.Range('D21').Value = '29.562.20'
$ws.Range('E21').Value = '  +0.31%  '
$ws.Range('D22').Value = '5.669'
$ws.Range('E22').Value = '  +0.89%  '
$ws.Range('D23').Value = '12.03'
$ws.Range('E23').Value = '  +2.76%  '
$ws.Range('D24').Value = '2.177'
$ws.Range('E24').Value = '  +0.76%  '
$ws.Range('D25').Value = '2.178.18'
$ws.Range('E25').Value = '  +2.70%  '
$ws.Range('D26').Value = '155.89'
$ws.Range('E26').Value = '  -0.24%  '
$ws.Range('D27').Value = '6.428'
$ws.Range('E27').Value = '  -0.54%  '
$ws.Range('D28').Value = '20.05'
$ws.Range('E28').Value = '  +0.09%  '
$ws.Range('D29').Value = '2.094'
$ws.Range('E29').Value = '  -0.83%  '
$ws.Range('D30').Value = '120.61'
$ws.Range('E30').Value = '  +0.50%  '
$ws.Range('D31').Value = '1.013'
$ws.Range('E31').Value = '  -1.18%  '
$ws.Range('D32').Value = '0.09583'
$ws.Range('E32').Value = '  +0.63%  '
$ws.Range('D33').Value = '5.596'
$ws.Range('E33').Value = '  +1.47%  '
$ws.Range('D34').Value = '3.566'
$ws.Range('E34').Value = '  +0.14%  '
$ws.Range('D35').Value = '1.384'
$ws.Range('E35').Value = '  -0.56%  '
$ws.Range('D36').Value = '0.06346'
$ws.Range('E36').Value = '  +4.04%  '
$ws.Range('E37').Value = '  +0.82%  '
$ws.Range('D38').Value = '1.195'
$ws.Range('E38').Value = '  +1.63%  '
$ws.Range('E39').Value = '  +0.03%  '
$ws.Range('D40').Value = '10.71'
$ws.Range('E40').Value = '  -1.23%  '
$ws.Range('D41').Value = '1.012'
$ws.Range('E41').Value = '  +0.66%  '
$ws.Range('D42').Value = '7.893'
$ws.Range('E42').Value = '  -1.23%  '
$ws.Range('D43').Value = '0.1846'
$ws.Range('E43').Value = '  -0.44%  '
$ws.Range('D44').Value = '2.481'
$ws.Range('E44').Value = '  +3.69%  '
$ws.Range('D45').Value = '1.280'
$ws.Range('E45').Value = '  +0.32%  '
$ws.Range('D46').Value = '12.44'
$ws.Range('E46').Value = '  -0.84%  '
$ws.Range('D47').Value = '0.07487'
$ws.Range('E47').Value = '  -1.45%  '
$ws.Range('D48').Value = '0.5551'
$ws.Range('E48').Value = '  -0.23%  '
$ws.Range('E49').Value = '  +1.62%  '
$ws.Range('D50').Value = '118.21'
$ws.Range('E50').Value = '  +1.58%  '
$ws.Range('D51').Value = '2.438'
$ws.Range('E51').Value = '  +1.46%  '

# Remove the temporary Text formatting so the cells end up with no
# explicit style, matching the original (unstyled) D/E data cells.
foreach ($ref in $targetRefs) {
    $ws.Range($ref).ClearFormats()
}
